# Update "想去人数" (want-to-go count) figures in column F that changed
# due to the data refresh captured at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rId1 / sheet1.xml
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 587
$ws1.Range("F8").Value  = 104
$ws1.Range("F9").Value  = 8670
$ws1.Range("F12").Value = 1140
$ws1.Range("F13").Value = 967
$ws1.Range("F14").Value = 105
$ws1.Range("F18").Value = 247
$ws1.Range("F19").Value = 65
$ws1.Range("F21").Value = 1010

# Sheet "全部类型" (All Types) - rId4 / sheet4.xml
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 587
$ws4.Range("F10").Value = 104
$ws4.Range("F11").Value = 8670
$ws4.Range("F14").Value = 1140
$ws4.Range("F15").Value = 967
$ws4.Range("F16").Value = 105
$ws4.Range("F20").Value = 247
$ws4.Range("F21").Value = 65
$ws4.Range("F23").Value = 1010
